$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value2 = 46875.59
$ws.Range("I28").Value2 = 60578.53
$ws.Range("K28").Value2 = 60578.53
$ws.Range("M28").Value2 = -60093.53

# Row 32
$ws.Range("H32").Value2 = 1592.4706
$ws.Range("J32").Value2 = 1663.1428
$ws.Range("L32").Value2 = 1663.1428
$ws.Range("N32").Value2 = -2315.1428

# Row 40
$ws.Range("H40").Value2 = 8528
$ws.Range("I40").Value2 = 7083.75
$ws.Range("J40").Value2 = 11416.5
$ws.Range("K40").Value2 = 7083.75
$ws.Range("L40").Value2 = 11416.5
$ws.Range("M40").Value2 = -6908.75
$ws.Range("N40").Value2 = -11766.5

# Row 43
$ws.Range("H43").Value2 = 1881.3704
$ws.Range("I43").Value2 = 1878.8422
$ws.Range("K43").Value2 = 1878.8422
$ws.Range("M43").Value2 = -1809.8422

# Row 76
$ws.Range("H76").Value2 = 76929830
$ws.Range("I76").Value2 = 7160
$ws.Range("J76").Value2 = 125006500
$ws.Range("K76").Value2 = 7160
$ws.Range("L76").Value2 = 125006500
$ws.Range("M76").Value2 = -6845
$ws.Range("N76").Value2 = -125007130

# Row 79
$ws.Range("H79").Value2 = 76929830
$ws.Range("I79").Value2 = 7160
$ws.Range("J79").Value2 = 125006500
$ws.Range("K79").Value2 = 7160
$ws.Range("L79").Value2 = 125006500
$ws.Range("M79").Value2 = -6068
$ws.Range("N79").Value2 = -125008684

# Row 98
$ws.Range("H98").Value2 = 708.86664
$ws.Range("I98").Value2 = 560.89655
$ws.Range("K98").Value2 = 560.89655
$ws.Range("M98").Value2 = 937.10345

# Row 100
$ws.Range("H100").Value2 = 5491
$ws.Range("I100").Value2 = 4815.3335
$ws.Range("K100").Value2 = 4815.3335
$ws.Range("M100").Value2 = -4274.3335

# Row 106
$ws.Range("H106").Value2 = 3729.652
$ws.Range("I106").Value2 = 3561.2856
$ws.Range("K106").Value2 = 3561.2856
$ws.Range("M106").Value2 = -2930.2856

# Row 122
$ws.Range("H122").Value2 = 708.86664
$ws.Range("I122").Value2 = 560.89655
$ws.Range("K122").Value2 = 1682.68965
$ws.Range("M122").Value2 = 767.3103499999997

# Row 132
$ws.Range("H132").Value2 = 2364.3784
$ws.Range("I132").Value2 = 2337.7354
$ws.Range("K132").Value2 = 7013.206200000001
$ws.Range("M132").Value2 = -4483.206200000001

# Row 137
$ws.Range("H137").Value2 = 1839.7059
$ws.Range("J137").Value2 = 2322.9412
$ws.Range("L137").Value2 = 6968.823600000001
$ws.Range("N137").Value2 = -12068.8236

# Row 138
$ws.Range("H138").Value2 = 5863.2256
$ws.Range("I138").Value2 = 2935.3333
$ws.Range("J138").Value2 = 7358.3193
$ws.Range("K138").Value2 = 8805.999899999999
$ws.Range("L138").Value2 = 22074.9579
$ws.Range("M138").Value2 = -3665.999899999999
$ws.Range("N138").Value2 = -32354.9579

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 12190.895
$ws.Range("I32").Value2 = 12424
$ws.Range("J32").Value2 = 7995
$ws.Range("K32").Value2 = 12424
$ws.Range("L32").Value2 = 7995
$ws.Range("M32").Value2 = -12137
$ws.Range("N32").Value2 = -8569

# Row 45
$ws.Range("H45").Value2 = 2058.56
$ws.Range("I45").Value2 = 1701.4706
$ws.Range("K45").Value2 = 1701.4706
$ws.Range("M45").Value2 = -1324.4706

# Row 61
$ws.Range("H61").Value2 = 4869.846
$ws.Range("I61").Value2 = 4525.6665
$ws.Range("J61").Value2 = 9000
$ws.Range("K61").Value2 = 4525.6665
$ws.Range("L61").Value2 = 9000
$ws.Range("M61").Value2 = -4313.6665
$ws.Range("N61").Value2 = -9424

# Row 102
$ws.Range("H102").Value2 = 2019.2
$ws.Range("I102").Value2 = 2048.0833
$ws.Range("J102").Value2 = 1903.6666
$ws.Range("K102").Value2 = 2048.0833
$ws.Range("L102").Value2 = 1903.6666
$ws.Range("M102").Value2 = -426.0832999999998
$ws.Range("N102").Value2 = -5147.6666

# Row 110
$ws.Range("H110").Value2 = 558929.5600000001
$ws.Range("I110").Value2 = 558929.5600000001
$ws.Range("K110").Value2 = 558929.5600000001
$ws.Range("M110").Value2 = -556884.5600000001

# Row 132
$ws.Range("H132").Value2 = 4582.6123
$ws.Range("I132").Value2 = 3452.2058
$ws.Range("J132").Value2 = 7144.8667
$ws.Range("K132").Value2 = 10356.6174
$ws.Range("L132").Value2 = 21434.6001
$ws.Range("M132").Value2 = -7826.617400000001
$ws.Range("N132").Value2 = -26494.6001

# Row 136
$ws.Range("H136").Value2 = 4869.846
$ws.Range("I136").Value2 = 4525.6665
$ws.Range("J136").Value2 = 9000
$ws.Range("K136").Value2 = 13576.9995
$ws.Range("L136").Value2 = 27000
$ws.Range("M136").Value2 = -11026.9995
$ws.Range("N136").Value2 = -32100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value2 = 151.75
$ws.Range("I80").Value2 = 115.8
$ws.Range("K80").Value2 = 115.8
$ws.Range("M80").Value2 = 882.2

# Row 83
$ws.Range("H83").Value2 = 151.75
$ws.Range("I83").Value2 = 115.8
$ws.Range("K83").Value2 = 579
$ws.Range("M83").Value2 = 4413

# Row 105
$ws.Range("H105").Value2 = 79755.92
$ws.Range("I105").Value2 = 79755.92
$ws.Range("K105").Value2 = 79755.92
$ws.Range("M105").Value2 = -78008.92

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value2 = 1039.5
$ws.Range("I105").Value2 = 976.3
$ws.Range("K105").Value2 = 976.3
$ws.Range("M105").Value2 = 770.7

# Row 122
$ws.Range("H122").Value2 = 1851.7894
$ws.Range("I122").Value2 = 1710.9412
$ws.Range("J122").Value2 = 3049
$ws.Range("K122").Value2 = 5132.8236
$ws.Range("L122").Value2 = 9147
$ws.Range("M122").Value2 = -2682.8236
$ws.Range("N122").Value2 = -14047

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value2 = 1917
$ws.Range("I23").Value2 = 1109.75
$ws.Range("J23").Value2 = 2275.7778
$ws.Range("K23").Value2 = 3329.25
$ws.Range("L23").Value2 = 6827.3334
$ws.Range("M23").Value2 = -3094.25
$ws.Range("N23").Value2 = -7297.3334

# Row 86
$ws.Range("H86").Value2 = 588.4
$ws.Range("I86").Value2 = 250
$ws.Range("J86").Value2 = 1096
$ws.Range("K86").Value2 = 750
$ws.Range("L86").Value2 = 3288
$ws.Range("M86").Value2 = 436
$ws.Range("N86").Value2 = -5660

# Row 89
$ws.Range("H89").Value2 = 588.4
$ws.Range("I89").Value2 = 250
$ws.Range("J89").Value2 = 1096
$ws.Range("K89").Value2 = 2250
$ws.Range("L89").Value2 = 9864
$ws.Range("M89").Value2 = 3678
$ws.Range("N89").Value2 = -21720

# Row 98
$ws.Range("H98").Value2 = 2948.4
$ws.Range("I98").Value2 = 3256
$ws.Range("J98").Value2 = 2816.5715
$ws.Range("K98").Value2 = 9768
$ws.Range("L98").Value2 = 8449.7145
$ws.Range("M98").Value2 = -8270
$ws.Range("N98").Value2 = -11445.7145

# Row 132
$ws.Range("H132").Value2 = 461480.97
$ws.Range("J132").Value2 = 628828.2
$ws.Range("L132").Value2 = 5659453.8
$ws.Range("N132").Value2 = -5664513.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value2 = 1258.4736
$ws.Range("I97").Value2 = 1307
$ws.Range("J97").Value2 = 999.6667
$ws.Range("K97").Value2 = 1307
$ws.Range("L97").Value2 = 999.6667
$ws.Range("M97").Value2 = -811
$ws.Range("N97").Value2 = -1991.6667

# Row 113
$ws.Range("H113").Value2 = 483678.38
$ws.Range("I113").Value2 = 910718.9399999999
$ws.Range("K113").Value2 = 910718.9399999999
$ws.Range("M113").Value2 = -908548.9399999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value2 = 499
$ws.Range("I22").Value2 = 499
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 499
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = -204
$ws.Range("N22").ClearContents()

# Row 27
$ws.Range("H27").Value2 = 499
$ws.Range("I27").Value2 = 499
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 499
$ws.Range("L27").Value2 = 0
$ws.Range("M27").Value2 = -392
$ws.Range("N27").ClearContents()

# Row 46
$ws.Range("H46").Value2 = 2499.3076
$ws.Range("I46").Value2 = 2644.3333
$ws.Range("J46").Value2 = 2173
$ws.Range("K46").Value2 = 2644.3333
$ws.Range("L46").Value2 = 2173
$ws.Range("M46").Value2 = -2456.3333
$ws.Range("N46").Value2 = -2549

# Row 55
$ws.Range("H55").Value2 = 899.6818
$ws.Range("I55").Value2 = 93.888885
$ws.Range("K55").Value2 = 93.888885
$ws.Range("M55").Value2 = 79.111115

# Row 100
$ws.Range("H100").Value2 = 2110.4546
$ws.Range("I100").Value2 = 1879.1177
$ws.Range("K100").Value2 = 1879.1177
$ws.Range("M100").Value2 = -1338.1177

# Row 110
$ws.Range("H110").Value2 = 61500
$ws.Range("J110").Value2 = 61500
$ws.Range("L110").Value2 = 61500
$ws.Range("N110").Value2 = -69680

# Row 122
$ws.Range("H122").Value2 = 5117.6855
$ws.Range("I122").Value2 = 4280.95
$ws.Range("J122").Value2 = 6233.3335
$ws.Range("K122").Value2 = 12842.85
$ws.Range("L122").Value2 = 18700.0005
$ws.Range("M122").Value2 = -10392.85
$ws.Range("N122").Value2 = -23600.0005

# Row 132
$ws.Range("H132").Value2 = 4097.6665
$ws.Range("J132").Value2 = 4333.3335
$ws.Range("L132").Value2 = 13000.0005
$ws.Range("N132").Value2 = -18060.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value2 = 5446.433
$ws.Range("I81").Value2 = 2278.3157
$ws.Range("J81").Value2 = 10918.637
$ws.Range("K81").Value2 = 4556.6314
$ws.Range("L81").Value2 = 21837.274
$ws.Range("M81").Value2 = -3495.6314
$ws.Range("N81").Value2 = -23959.274

# Row 84
$ws.Range("H84").Value2 = 5446.433
$ws.Range("I84").Value2 = 2278.3157
$ws.Range("J84").Value2 = 10918.637
$ws.Range("K84").Value2 = 22783.157
$ws.Range("L84").Value2 = 109186.37
$ws.Range("M84").Value2 = -17479.157
$ws.Range("N84").Value2 = -119794.37
